$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet/tab
$ws.Name = "Sheet1"

# Insert 4 new blank rows at the top; existing rows 1-15 become rows 5-19
$ws.Rows("1:4").Insert()

# Copy formatting (border/alignment/wrap) from the existing formatted row
# down onto the new header block (rows 1-3, columns A:D only) without
# introducing any new/unused cell styles.
$ws.Range("A5:D5").Copy()
$ws.Range("A1:D3").PasteSpecial(-4122)

# Fill in the new "master package" header block text
$ws.Range("A1").Value = "MASTER PACKAGE"

$ws.Range("A2").Value = "WesternGlove Centric8 PROD"
$ws.Range("B2").Value = "M12225BVS563:KONRAD"
$ws.Range("C2").Value = "DUP REVIEW"
$ws.Range("D2").Value = "Revised 1/9/25, 2:47 PM"

$ws.Range("A3").Value = "Evaluation"

# New trailing rows with additional table data
$ws.Range("A16").Value = "B170"
$ws.Range("B16").Value = "BELTLOOP WIDTH"
$ws.Range("D16").Value = "A6"
$ws.Range("E16").Value = "%"
$ws.Range("F16").Value = "%"
$ws.Range("H16").Value = "%"

$ws.Range("A17").Value = "B171"
$ws.Range("B17").Value = "BELTLOOP LENGTH (FINISHED)"
$ws.Range("F17").Value = "2%"
$ws.Range("H17").Value = "2%"

$ws.Range("A18").Value = "B172"
$ws.Range("B18").Value = "BELTLOOP LENGTH (TACK TO TACK)"
$ws.Range("F18").Value = "2%"
$ws.Range("H18").Value = "2%"

$ws.Range("A19").Value = "MBP1 - CONTEMPORARY"

# Adjust column widths (A, C, D) to match the new table layout.
# The Excel COM ColumnWidth property is quantized to on-screen pixel
# steps, so these inputs are chosen to land on the closest achievable
# stored width.
$ws.Columns("A").ColumnWidth = 26.0
$ws.Columns("C").ColumnWidth = 11.166666666666666
$ws.Columns("D").ColumnWidth = 21.166666666666668
